$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("var read across ")

# Copy the 2022 column (J) data and formatting into the new 2023 column (K)
$ws.Range("J1:J43").Copy()
$ws.Range("K1:K43").PasteSpecial(-4163)

# Update the new 2023 column's header value
$ws.Range("K1").Value = 2023
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").HorizontalAlignment = -4131

$ws.Range("K2").Select()
